$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (item 8): 2018-08-09, RFM95 868MHz, 250000, Thao, TaoBao ---
$ws.Range("B11").Value = 43321
$ws.Range("C11").Value = "RFM95 868MHz"
$ws.Range("D11").Value = 250000
$ws.Range("E11").Value = "Thao"
$ws.Range("F11").Value = "TaoBao"

# --- Row 12 (item 9): 2018-08-09, PCB of LCD_Button_Gateway_Driver & RPi3_LoRaWan_Gateway, 350000, Thao, JLCPCB ---
$ws.Range("B12").Value = 43321
$ws.Range("C12").Value = "PCB of LCD_Button_Gateway_Driver & RPi3_LoRaWan_Gateway"
$ws.Range("D12").Value = 350000
$ws.Range("E12").Value = "Thao"
$ws.Range("F12").Value = "JLCPCB"

# --- Row 13 (item 10): 2018-08-05, (item name left blank), 70400, Thao, Thegioiic, Bill number 34246 ---
$ws.Range("B13").Value = 43317
$ws.Range("D13").Value = 70400
$ws.Range("E13").Value = "Thao"
$ws.Range("F13").Value = "Thegioiic"
$ws.Range("G13").Value = "Bill number 34246"

# --- Row 14 (item 11): 2018-08-09, 2 x 18650 Battery 3v3 + ATMega328P-AU, 293000, Thao, Thegioiic, Bill number 34506 ---
# Note field (G14) is written before the item-name field (C14) to keep the new
# shared-string insertion order identical to the source edit.
$ws.Range("B14").Value = 43321
$ws.Range("D14").Value = 293000
$ws.Range("E14").Value = "Thao"
$ws.Range("F14").Value = "Thegioiic"
$ws.Range("G14").Value = "Bill number 34506"
$ws.Range("C14").Value = "2 x 18650 Battery 3v3 + ATMega328P-AU"

# The date values above were written as plain serials so no auto "guessed" date
# format gets synthesized as a brand-new style; now stamp them with the same
# date format (style) already used by the existing date column (B10) via a
# format-only paste, which reuses the existing style instead of creating one.
$ws.Range("B10").Copy()
$ws.Range("B11:B14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update selection / active cell to C13 (matches the saved sheetView selection)
$ws.Range("C13").Select()

# Widen column B. The target OOXML stored width is 33.7109375; this runtime's
# ColumnWidth setter quantizes to 1/6-character steps, so 33.7109375 itself is
# unreachable. 32.8 is the input that lands on the nearest achievable stored
# width (33.666666666666664).
$ws.Columns.Item(2).ColumnWidth = 32.8
